$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "ongules" (Réseau Ongulés sauvages) row (row 4)
$ws.Range("B4").Value = "Réseau Ongulés sauvages"
$ws.Range("D4").Value = "Le réseau a pour but de récolter les informations utiles pour suivre les 14 espèces d'ongulés sauvages présentes en France hexagonale"
$ws.Range("E4").Value = "L’ensemble des données récoltées permet de mesurer de façon régulière des variables biologiques d’intérêt, telles que les aires de présence, les tendances d’évolution ou l’état de santé des populations d’ongulés sauvages en France métropolitaine. En parallèle, d’autres données connexes, telles les prélèvements cynégétiques ou les modalités de gestion des populations, sont régulièrement enregistrées."
$ws.Range("G4").Value = "75,77,78,91,92,93,94,95"
$ws.Range("J4").Value = "1,2,3,4,5,6,7,8,9,10,11,12"
$ws.Range("L4").Value = "Animation nationale: XXXX`nAnimation régionale: Samuel DEMBSKI`nCorrespondants départementaux:`nPPC:`n77:`n78-95:`n91`nCourriel du réseau: reseau.ongules-sauvages@ofb.gouv.fr"
$ws.Range("M4").Value = "Fédérations de chasse"
$ws.Range("Z4").Value = "- Articles de recherche`n- Lettres d'informations`n- Bilan annuel`n- Cartes de densité"
$ws.Range("AA4").Value = "texte:Dataviz: Présence des ongulés sauvages en France;lien:https://professionnels.ofb.fr/fr/doc-dataviz/dataviz-presence-ongules-sauvages-en-France"
$ws.Range("AB4").Value = "texte:Fiches de synthèse des suivis;lien:https://professionnels.ofb.fr/fr/node/869"
$ws.Range("AC4").Value = "texte: Cartes de répartition;lien:https://carmen.carmencarto.fr/38/Ongules_sauvages.map#"
$ws.Range("AD4").Value = "texte:Lettre d'information;lien:https://professionnels.ofb.fr/fr/node/1281"
$ws.Range("AE4").Value = "texte: Page du réseau sur le portail technique;lien:https://professionnels.ofb.fr/node/1431"

# Setting the value on Z4 resets its cell style from the highlighted "diffusion"
# style back to default; restore it by copying the format from a neighbouring
# cell (U4) that already carries the same style.
$ws.Range("U4").Copy()
$ws.Range("Z4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 4 needs to grow to fit the new multi-line content, matching row 5's height
$ws.Rows.Item(4).RowHeight = 285

# Update the active selection in the frozen pane to reflect the new edited cell
$ws.Range("J4").Select()
